$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new event row (Wrathborne Invasion 12PM) under the existing data
$ws.Range("A12").Value = "Wrathborne Invasion 12PM"

# Move the active selection to match the post-edit state
$ws.Range("F14").Select()
